$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Save" header in H1, matching style/formatting of existing headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H7 with 0 values
$ws.Range("H2:H7").Value = 0
